$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - values are stored as text, matching the
# original inline-string cell type, so force a Text number format before
# assigning the numeric-looking string value.
$priceUpdates = @{
    2 = "302.13"
    3 = "35.16"
    4 = "5.035"
    5 = "0.07904"
    6 = "1.946"
    7 = "7.747"
    8 = "4.027"
    9 = "2.871"
    10 = "0.9228"
    11 = "0.1194"
    12 = "0.1837"
    13 = "0.09315"
    14 = "0.03537"
    15 = "0.09868"
    16 = "0.001388"
    17 = "0.005930"
    18 = "3.493"
    19 = "0.3442"
    20 = "0.1309"
    22 = "0.2399"
    23 = "0.04505"
    39 = "0.01896"
    40 = "0.04709"
    41 = "0.007566"
    43 = "0.1323"
    45 = "0.01119"
    46 = "0.00006025"
    50 = "0.00002100"
    51 = "0.0002000"
}

# Volume(1h) (column E) updates - percentages stored as text.
$volumeUpdates = @{
    2 = "-6.26%"
    3 = "-3.14%"
    4 = "-1.64%"
    5 = "-2.29%"
    6 = "-9.73%"
    7 = "-3.46%"
    8 = "-2.81%"
    9 = "2.51%"
    10 = "-0.59%"
    11 = "18.20%"
    12 = "-2.47%"
    13 = "1.76%"
    14 = "-1.65%"
    15 = "-0.73%"
    16 = "-3.28%"
    17 = "4.30%"
    18 = "1.11%"
    19 = "2.09%"
    20 = "-1.71%"
    21 = "-0.63%"
    22 = "8.86%"
    23 = "-2.15%"
    24 = "-2.45%"
    25 = "-3.35%"
    26 = "-3.89%"
    27 = "-6.96%"
    39 = "-6.73%"
    40 = "-5.59%"
    41 = "-2.84%"
    42 = "22.04%"
    43 = "-5.46%"
    44 = "1.37%"
    45 = "-7.86%"
    46 = "-6.73%"
    47 = "-0.11%"
    49 = "-31.45%"
    50 = "-0.11%"
    51 = "-0.11%"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $volumeUpdates[$row]
}

# Hora (column G) - every data row (2-51) moves from "4" to "5".
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = "5"
}
